# Fruta / hortaliza, semanal
# Insert two new weekly report rows before the current row 417, shifting the
# existing rows 417:442 down to 419:444, and populate the two new rows with
# the latest week's Pimiento data for "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 417 (each Insert() pushes rows 417+ down by one).
$ws.Rows.Item(417).Insert()
$ws.Rows.Item(417).Insert()

# New row 417: Cuatro cascos rojo
$ws.Cells.Item(417, 1).Value = 5
$ws.Cells.Item(417, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(417, 3).Value = "Maule"
$ws.Cells.Item(417, 4).Value2 = 44610
$ws.Cells.Item(417, 5).Value = 7
$ws.Cells.Item(417, 6).Value = 100112002
$ws.Cells.Item(417, 7).Value = "Pimiento"
$ws.Cells.Item(417, 8).Value = "Cuatro cascos rojo"
$ws.Cells.Item(417, 9).Value = "Primera"
$ws.Cells.Item(417, 10).Value = 200
$ws.Cells.Item(417, 11).Value = 10000
$ws.Cells.Item(417, 12).Value = 10000
$ws.Cells.Item(417, 13).Value = 10000
$ws.Cells.Item(417, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(417, 15).Value = "Región del Maule"
$ws.Cells.Item(417, 16).Value = 667
$ws.Cells.Item(417, 17).Value = 15
$ws.Cells.Item(417, 18).Value = "Hortaliza"

# New row 418: Cuatro cascos verde
$ws.Cells.Item(418, 1).Value = 5
$ws.Cells.Item(418, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(418, 3).Value = "Maule"
$ws.Cells.Item(418, 4).Value2 = 44610
$ws.Cells.Item(418, 5).Value = 7
$ws.Cells.Item(418, 6).Value = 100112002
$ws.Cells.Item(418, 7).Value = "Pimiento"
$ws.Cells.Item(418, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(418, 9).Value = "Primera"
$ws.Cells.Item(418, 10).Value = 300
$ws.Cells.Item(418, 11).Value = 6000
$ws.Cells.Item(418, 12).Value = 6000
$ws.Cells.Item(418, 13).Value = 6000
$ws.Cells.Item(418, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(418, 15).Value = "Región del Maule"
$ws.Cells.Item(418, 16).Value = 400
$ws.Cells.Item(418, 17).Value = 15
$ws.Cells.Item(418, 18).Value = "Hortaliza"
